$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,('Agency','Program','Checked')
    ,('WRM Agency','National Water Resource Management Sector Framework','Y')
    ,('WRM Agency','Overarching Flood Risk Management Framework','Y')
    ,('WRM Agency','Overarching National Drought Risk Management Framework','Y')
    ,('WRM Agency','Local Government','')
    ,('WRM Agency','Public Participation & Stakeholder Engagement','')
    ,('WRM Agency','Social Inclusion','')
    ,('WRM Agency','Education & Risk Communication','')
    ,('WRM Agency','Scientific Collaboration ','')
    ,('WRM Agency','Open Data','')
    ,('WRM Agency','Integrated River Basin Planning','Y')
    ,('WRM Agency','Urban Water Supply Planning','Y')
    ,('WRM Agency','Irrigation Water Supply Planning','Y')
    ,('WRM Agency','Water Resources Investment Policy','Y')
    ,('WRM Agency','Dam Safety','Y')
    ,('WRM Agency','Flood Infrastructure Safety','Y')
    ,('WRM Agency','Flexible Water Allocation','Y')
    ,('WRM Agency','Conjunctive Groundwater Management','Y')
    ,('WRM Agency','Floodplain Mapping','Y')
    ,('WRM Agency','Local Flood Risk Mitigation Planning','y')
    ,('WRM Agency','Drought Monitoring Program','Y')
    ,('WRM Agency','WRM Drought Response','Y')
    ,('WRM Agency','Flood Forecasting and Warning','Y')
    ,('DRM Agency','National Disaster Risk Management Sector Framework','Y')
    ,('DRM Agency','Overarching Flood Risk Management Framework','Y')
    ,('DRM Agency','Local Government','')
    ,('DRM Agency','Public Participation & Stakeholder Engagement','')
    ,('DRM Agency','Social Inclusion','')
    ,('DRM Agency','Education & Risk Communication','')
    ,('DRM Agency','Scientific Collaboration ','')
    ,('DRM Agency','Open Data','')
    ,('DRM Agency','Local Flood Mitigation Planning','y')
    ,('DRM Agency','Floodplain Mapping','Y')
    ,('DRM Agency','Floodplain Regulation','Y')
    ,('DRM Agency','Local Flood Mitigation Planning','')
    ,('DRM Agency','Flood Forecasting and Warning','Y')
    ,('DRM Agency','Flood Emergency Preparedness, Response, and Relief','')
    ,('DRM Agency','Flood Disaster Recovery','Y')
    ,('Hydro-Met','Local Government','')
    ,('Hydro-Met','Public Participation & Stakeholder Engagement','')
    ,('Hydro-Met','Social Inclusion','')
    ,('Hydro-Met','Education & Risk Communication','')
    ,('Hydro-Met','Scientific Collaboration ','')
    ,('Hydro-Met','Open Data','')
    ,('Hydro-Met','National Framework for NMS/NHS Services','Y')
    ,('Hydro-Met','Co-Production of Services','Y')
    ,('Hydro-Met','Flood Forecasting and Warning','Y')
    ,('Agriculture','Overarching National Drought Risk Management Framework','Y')
    ,('Agriculture','Local Government','')
    ,('Agriculture','Public Participation & Stakeholder Engagement','')
    ,('Agriculture','Social Inclusion','')
    ,('Agriculture','Education & Risk Communication','')
    ,('Agriculture','Scientific Collaboration ','')
    ,('Agriculture','Open Data','')
    ,('Agriculture','Irrigation Water Supply Planning','Y')
    ,('Agriculture','Agriculture Policies and Climate Smart Agriculture','Y')
    ,('Agriculture','Watershed Management','')
    ,('Agriculture','Agriculture Drought Response','Y')
    ,('Agriculture','Social Protection Drought Response','Y')
    ,('Natural Resources Management ','Local Government','')
    ,('Natural Resources Management ','Public Participation & Stakeholder Engagement','')
    ,('Natural Resources Management ','Social Inclusion','')
    ,('Natural Resources Management ','Education & Risk Communication','')
    ,('Natural Resources Management ','Scientific Collaboration ','')
    ,('Natural Resources Management ','Open Data','')
    ,('Natural Resources Management ','Coastal Zone Management Planning','Y')
    ,('Natural Resources Management ','Forest Management','Y')
    ,('Natural Resources Management ','Wetlands Management','Y')
    ,('Natural Resources Management ','Watershed Management','Y')
    ,('Social Protection ','Local Government','')
    ,('Social Protection ','Public Participation & Stakeholder Engagement','')
    ,('Social Protection ','Social Inclusion','')
    ,('Social Protection ','Education & Risk Communication','')
    ,('Social Protection ','Scientific Collaboration ','')
    ,('Social Protection ','Open Data','')
    ,('Social Protection ','Social Protection Drought Response','Y')
    ,('Finance Agency','Local Government','')
    ,('Finance Agency','Public Participation & Stakeholder Engagement','')
    ,('Finance Agency','Social Inclusion','')
    ,('Finance Agency','Education & Risk Communication','')
    ,('Finance Agency','Scientific Collaboration ','')
    ,('Finance Agency','Open Data','')
    ,('Finance Agency','Disaster Risk Financing Instruments','')
    ,('Finance Agency','Disaster Risk Financing National Sector Framework','')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    if ($row[2] -ne "") {
        $ws.Cells.Item($r, 3).Value = $row[2]
    } else {
        $ws.Cells.Item($r, 3).ClearContents()
    }
}

# Apply white fill to B37 (DRM Agency / Flood Emergency Preparedness, Response, and Relief)
$ws.Range("B37").Interior.ThemeColor = 2
